$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Value)
    $Cell.Value = "'" + $Value
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) '28.614.84'
$ws.Cells.Item(2, 5).Value = '  +1.08%  '
Set-TextValue $ws.Cells.Item(3, 4) '1.558.10'
$ws.Cells.Item(3, 5).Value = '  -1.04%  '
Set-TextValue $ws.Cells.Item(4, 4) '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.28%  '
Set-TextValue $ws.Cells.Item(5, 4) '210.82'
$ws.Cells.Item(5, 5).Value = '  -0.52%  '
Set-TextValue $ws.Cells.Item(6, 4) '0.486'
$ws.Cells.Item(6, 5).Value = '  -0.69%  '
$ws.Cells.Item(7, 5).Value = '  -0.33%  '
Set-TextValue $ws.Cells.Item(8, 4) '24.57'
$ws.Cells.Item(8, 5).Value = '  +3.26%  '
$ws.Cells.Item(9, 5).Value = '  -0.25%  '
$ws.Cells.Item(10, 5).Value = '  -0.53%  '
Set-TextValue $ws.Cells.Item(11, 4) '0.0894'
$ws.Cells.Item(11, 5).Value = '  -0.05%  '
Set-TextValue $ws.Cells.Item(12, 4) '1.779.71'
$ws.Cells.Item(12, 5).Value = '  -1.14%  '
Set-TextValue $ws.Cells.Item(13, 4) '1.559.26'
$ws.Cells.Item(13, 5).Value = '  -1.45%  '
Set-TextValue $ws.Cells.Item(14, 4) '28.617.88'
$ws.Cells.Item(14, 5).Value = '  +0.97%  '
Set-TextValue $ws.Cells.Item(15, 4) '0.514'
$ws.Cells.Item(15, 5).Value = '  -0.42%  '
Set-TextValue $ws.Cells.Item(16, 4) '3.64'
$ws.Cells.Item(16, 5).Value = '  -1.29%  '
Set-TextValue $ws.Cells.Item(17, 4) '61.41'
$ws.Cells.Item(17, 5).Value = '  -0.28%  '
Set-TextValue $ws.Cells.Item(18, 4) '229.83'
$ws.Cells.Item(18, 5).Value = '  -0.11%  '
Set-TextValue $ws.Cells.Item(19, 4) '7.38'
$ws.Cells.Item(19, 5).Value = '  -0.47%  '
Set-TextValue $ws.Cells.Item(20, 4) '0.0₃0672'
$ws.Cells.Item(20, 5).Value = '  -1.81%  '
Set-TextValue $ws.Cells.Item(21, 4) '0.999'
$ws.Cells.Item(21, 5).Value = '  -0.25%  '
Set-TextValue $ws.Cells.Item(22, 4) '3.92'
$ws.Cells.Item(22, 5).Value = '  -0.88%  '
Set-TextValue $ws.Cells.Item(23, 4) '8.97'
$ws.Cells.Item(23, 5).Value = '  -0.70%  '
Set-TextValue $ws.Cells.Item(24, 4) '2.08'
$ws.Cells.Item(24, 5).Value = '  +1.26%  '
Set-TextValue $ws.Cells.Item(25, 4) '151.18'
$ws.Cells.Item(25, 5).Value = '  -0.31%  '
Set-TextValue $ws.Cells.Item(26, 4) '14.77'
$ws.Cells.Item(26, 5).Value = '  -1.10%  '
$ws.Cells.Item(27, 5).Value = '  -0.32%  '
$ws.Cells.Item(28, 5).Value = '  -0.20%  '
$ws.Cells.Item(29, 5).Value = '  -2.00%  '
Set-TextValue $ws.Cells.Item(30, 4) '0.0460'
$ws.Cells.Item(30, 5).Value = '  -4.24%  '
$ws.Cells.Item(31, 5).Value = '  -1.47%  '
Set-TextValue $ws.Cells.Item(32, 4) '3.17'
$ws.Cells.Item(32, 5).Value = '  -0.85%  '
Set-TextValue $ws.Cells.Item(33, 4) '1.391.18'
$ws.Cells.Item(33, 5).Value = '  +0.55%  '
$ws.Cells.Item(34, 5).Value = '  -2.32%  '
$ws.Cells.Item(35, 5).Value = '  -2.83%  '
$ws.Cells.Item(36, 5).Value = '  -1.75%  '
$ws.Cells.Item(37, 5).Value = '  +0.37%  '
Set-TextValue $ws.Cells.Item(38, 4) '2.28'
$ws.Cells.Item(38, 5).Value = '  -3.58%  '
$ws.Cells.Item(39, 5).Value = '  -0.66%  '
Set-TextValue $ws.Cells.Item(40, 4) '1.95'
$ws.Cells.Item(40, 5).Value = '  +3.60%  '
Set-TextValue $ws.Cells.Item(41, 4) '0.517'
$ws.Cells.Item(41, 5).Value = '  -0.15%  '
$ws.Cells.Item(42, 5).Value = '  -0.23%  '
Set-TextValue $ws.Cells.Item(43, 4) '0.776'
$ws.Cells.Item(43, 5).Value = '  -1.20%  '
Set-TextValue $ws.Cells.Item(44, 4) '0.0464'
$ws.Cells.Item(44, 5).Value = '  +0.15%  '
Set-TextValue $ws.Cells.Item(45, 4) '64.02'
$ws.Cells.Item(45, 5).Value = '  +2.77%  '
Set-TextValue $ws.Cells.Item(46, 4) '5.29'
$ws.Cells.Item(46, 5).Value = '  -1.59%  '
Set-TextValue $ws.Cells.Item(47, 4) '1.693.86'
$ws.Cells.Item(47, 5).Value = '  -1.05%  '
Set-TextValue $ws.Cells.Item(48, 4) '0.870'
$ws.Cells.Item(49, 2).Value = 'Quant'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Cells.Item(49, 4) '85.17'
$ws.Cells.Item(49, 5).Value = '  -0.42%  '
$ws.Cells.Item(50, 2).Value = 'BitcoinSV'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue $ws.Cells.Item(50, 4) '43.38'
$ws.Cells.Item(50, 5).Value = '  +4.86%  '
Set-TextValue $ws.Cells.Item(51, 4) '0.0₆0102'
$ws.Cells.Item(51, 5).Value = '  -0.13%  '
